# Apply cryptos list update (GitHub Actions refresh) to Sheet1
function Set-CellText($range, $text) {
    $c = $ws.Range($range)
    # Leading apostrophe forces Excel to store the value as text even
    # when it looks like a number (e.g. "604.06", "2.71"), matching the
    # source data where every cell in these columns is a string.
    $c.Value = "'" + $text
    # Reset formatting so the text-entry quote-prefix does not linger
    # as a visible style on the cell.
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText "D2" '64.058.06'
Set-CellText "D3" '3.127.73'
Set-CellText "E3" '  -2.68%  '
Set-CellText "D5" '604.06'
Set-CellText "E5" '  -0.54%  '
Set-CellText "D6" '147.00'
Set-CellText "E6" '  -5.31%  '
Set-CellText "E7" '  +0.07%  '
Set-CellText "D8" '3.123.83'
Set-CellText "E8" '  -2.71%  '
Set-CellText "D9" '0.525'
Set-CellText "E9" '  -3.91%  '
Set-CellText "D10" '0.151'
Set-CellText "E10" '  -5.63%  '
Set-CellText "D11" '5.54'
Set-CellText "E11" '  -2.58%  '
Set-CellText "D12" '0.472'
Set-CellText "E12" '  -5.58%  '
Set-CellText "D13" '0.0000256'
Set-CellText "E13" '  -3.97%  '
Set-CellText "D14" '36.25'
Set-CellText "E14" '  -5.13%  '
Set-CellText "D15" '3.648.83'
Set-CellText "E15" '  -2.54%  '
Set-CellText "D16" '64.093.64'
Set-CellText "E16" '  -3.46%  '
Set-CellText "D17" '3.148.48'
Set-CellText "E17" '  -2.15%  '
Set-CellText "E18" '  -0.06%  '
Set-CellText "D19" '6.90'
Set-CellText "E19" '  -4.65%  '
Set-CellText "D20" '477.68'
Set-CellText "E20" '  -5.54%  '
Set-CellText "D21" '14.46'
Set-CellText "E21" '  -4.99%  '
Set-CellText "D22" '0.703'
Set-CellText "E22" '  -3.27%  '
Set-CellText "D23" '7.64'
Set-CellText "E23" '  -4.25%  '
Set-CellText "D24" '13.64'
Set-CellText "E24" '  -5.78%  '
Set-CellText "D25" '83.26'
Set-CellText "E25" '  -1.89%  '
Set-CellText "D27" '2.91'
Set-CellText "E27" '  -2.79%  '
Set-CellText "D28" '8.43'
Set-CellText "E28" '  -5.85%  '
Set-CellText "D29" '2.22'
Set-CellText "E29" '  -5.31%  '
Set-CellText "E30" '  -16.59%  '
Set-CellText "D31" '6.81'
Set-CellText "E31" '  -1.39%  '
Set-CellText "B32" 'FirstDigitalUSD'
Set-CellText "C32" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText "D32" '1.00'
Set-CellText "E32" '  -0.01%  '
Set-CellText "B33" 'Stacks'
Set-CellText "C33" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText "D33" '2.71'
Set-CellText "E33" '  -6.10%  '
Set-CellText "D34" '26.45'
Set-CellText "E34" '  -6.07%  '
Set-CellText "D35" '1.10'
Set-CellText "E35" '  -5.32%  '
Set-CellText "D36" '6.03'
Set-CellText "E36" '  -5.39%  '
Set-CellText "D37" '54.44'
Set-CellText "E37" '  -1.70%  '
Set-CellText "D38" '3.10'
Set-CellText "E38" '  +3.40%  '
Set-CellText "D39" '0.0₃0731'
Set-CellText "E39" '  -4.89%  '
Set-CellText "D40" '447.91'
Set-CellText "E40" '  -10.07%  '
Set-CellText "D41" '0.0396'
Set-CellText "E41" '  -5.51%  '
Set-CellText "E42" '  -5.66%  '
Set-CellText "D43" '8.37'
Set-CellText "E43" '  -3.72%  '
Set-CellText "D44" '2.853.94'
Set-CellText "E44" '  -2.10%  '
Set-CellText "D45" '0.268'
Set-CellText "E45" '  -8.50%  '
Set-CellText "D46" '2.26'
Set-CellText "E46" '  -6.80%  '
Set-CellText "D47" '26.36'
Set-CellText "E47" '  -5.57%  '
Set-CellText "E48" '  -0.02%  '
Set-CellText "E49" '  -2.68%  '
Set-CellText "D50" '2.29'
Set-CellText "E50" '  -3.75%  '
Set-CellText "D51" '119.22'
Set-CellText "E51" '  -2.00%  '
